$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "consumer/typeset/runme_large.sh"
$ws.Range("B9").Value = 0.05
$ws.Range("C9").Value = 0.05
$ws.Range("D9").Value = 0
